$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Match the saved window geometry from the author's session
$win = $wb.Windows.Item(1)
$win.Left = 4185
$win.Top = 1365
$win.Width = 15375
$win.Height = 7875

# Add four new customer test rows to the existing AddCustomerTest sheet
$ws1.Range("A3").Value = "rahul"
$ws1.Range("B3").Value = "sahoo"
$ws1.Range("C3").Value = 751005
$ws1.Range("D3").Value = "Customer added successfully"

$ws1.Range("A4").Value = "tintun"
$ws1.Range("B4").Value = "sahoo"
$ws1.Range("C4").Value = 751005
$ws1.Range("D4").Value = "Customer added successfully"

$ws1.Range("A5").Value = "lukeworm"
$ws1.Range("B5").Value = "sahoo"
$ws1.Range("C5").Value = 751009
$ws1.Range("D5").Value = "Customer added successfully"

$ws1.Range("A6").Value = "smart"
$ws1.Range("B6").Value = "baby"
$ws1.Range("C6").Value = 751005
$ws1.Range("D6").Value = "Customer added successfully"

# Best-fit width for the new first column (renders as width="10")
$ws1.Columns.Item(1).ColumnWidth = 9.140625

# Move the selection on the first sheet
$ws1.Range("C4").Select()

# Add the new OpenAccountTest sheet right after AddCustomerTest
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"

$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "smruti sahoo"
$ws2.Range("B2").Value = "Rupee"

$ws2.Columns.Item(1).ColumnWidth = 11.5

$ws2.Range("A2").Select()

Write-Host "done"
